$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 16:22"

# Row 4
$ws.Range("B4").Value = 312152
$ws.Range("C4").Value = 795
$ws.Range("E4").Value = 288858
$ws.Range("G4").Value = 14
$ws.Range("H4").Value = 8466

# Row 17
$ws.Range("B17").Value = 11920
$ws.Range("C17").Value = 139
$ws.Range("E17").Value = 8718

# Row 21
$ws.Range("E21").Value = 7494
$ws.Range("G21").Value = 3
$ws.Range("H21").Value = 47

# Row 23
$ws.Range("F23").Value = 91

# Row 50
$ws.Range("F50").Value = 1

# Row 51
$ws.Range("A51").Value = "Islandia"
$ws.Range("B51").Value = 1486
$ws.Range("C51").Value = 69
$ws.Range("D51").Value = 428
$ws.Range("E51").Value = 1054
$ws.Range("F51").Value = 11
$ws.Range("H51").Value = 4

# Row 52
$ws.Range("A52").Value = "Argentina"
$ws.Range("B52").Value = 1451
$ws.Range("D52").Value = 279
$ws.Range("E52").Value = 1128
$ws.Range("F52").Value = 86
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 44

# Row 72
$ws.Range("B72").Value = 654
$ws.Range("C72").Value = 30
$ws.Range("E72").Value = 603

# Row 90
$ws.Range("A90").Value = "Reunion"
$ws.Range("B90").Value = 344
$ws.Range("C90").Value = 10
$ws.Range("D90").Value = 40
$ws.Range("E90").Value = 304
$ws.Range("F90").Value = 4
$ws.Range("H90").Value = 0

# Row 91
$ws.Range("A91").Value = "Afganistan"
$ws.Range("B91").Value = 337
$ws.Range("C91").Value = 28
$ws.Range("D91").Value = 15
$ws.Range("E91").Value = 315
$ws.Range("F91").Value = 0
$ws.Range("H91").Value = 7

# Row 117
$ws.Range("A117").Value = "Kenia"
$ws.Range("B117").Value = 142
$ws.Range("C117").Value = 16
$ws.Range("D117").Value = 4
$ws.Range("E117").Value = 134
$ws.Range("F117").Value = 2
$ws.Range("H117").Value = 4

# Row 118
$ws.Range("A118").Value = "Brunei"
$ws.Range("B118").Value = 135
$ws.Range("D118").Value = 73
$ws.Range("E118").Value = 61
$ws.Range("H118").Value = 1

# Row 119
$ws.Range("A119").Value = "Mayotte"
$ws.Range("D119").Value = 14
$ws.Range("E119").Value = 118
$ws.Range("F119").Value = 3
$ws.Range("H119").Value = 2

# Row 120
$ws.Range("A120").Value = "Guadalupe"
$ws.Range("B120").Value = 134
$ws.Range("C120").Value = 0
$ws.Range("D120").Value = 24
$ws.Range("E120").Value = 103
$ws.Range("F120").Value = 14
$ws.Range("H120").Value = 7

# Row 121
$ws.Range("A121").Value = "Isla de Man"
$ws.Range("B121").Value = 127
$ws.Range("C121").Value = 1
$ws.Range("D121").Value = 0
$ws.Range("E121").Value = 126
$ws.Range("F121").Value = 0
$ws.Range("H121").Value = 1

# Row 124
$ws.Range("A124").Value = "Trinidad yTobago"
$ws.Range("C124").Value = 1
$ws.Range("D124").Value = 1
$ws.Range("E124").Value = 96
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 1
$ws.Range("H124").Value = 7

# Row 125
$ws.Range("A125").Value = "Paraguay"
$ws.Range("B125").Value = 104
$ws.Range("C125").Value = 8
$ws.Range("D125").Value = 12
$ws.Range("E125").Value = 89
$ws.Range("F125").Value = 2
$ws.Range("H125").Value = 3

# Row 130
$ws.Range("B130").Value = 72
$ws.Range("C130").Value = 2
$ws.Range("D130").Value = 2

# Row 142
$ws.Range("A142").Value = "Togo"
$ws.Range("B142").Value = 44
$ws.Range("C142").Value = 3
$ws.Range("D142").Value = 20
$ws.Range("E142").Value = 21
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 3

# Row 143
$ws.Range("A143").Value = "Etiopia"
$ws.Range("B143").Value = 43
$ws.Range("C143").Value = 5
$ws.Range("D143").Value = 4
$ws.Range("E143").Value = 38
$ws.Range("F143").Value = 1
$ws.Range("G143").Value = 1
$ws.Range("H143").Value = 1

# Row 144
$ws.Range("A144").Value = "Mali"
$ws.Range("D144").Value = 1
$ws.Range("E144").Value = 37
